$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "unknown" peaks (and one mislabeled "me-17:0" / duplicated "C23:0 (IS)")
# with their identified compound names. Order matters: it reproduces the order in
# which the shared-string table grows with newly introduced labels.
$ws.Range("A52").Value = "C18:1 cis"
$ws.Range("A53").Value = "C18:1 trans"
$ws.Range("A57").Value = "C18:0 -me"
$ws.Range("A39").Value = "C17:0 -me"
$ws.Range("A59").Value = "cyc-19"
$ws.Range("A65").Value = "C20:0"
$ws.Range("A74").Value = "C22:0"
$ws.Range("A78").Value = "C23:0 (IS)"
$ws.Range("A82").Value = "alkane"
$ws.Range("A83").Value = "C24:0"

# Move the active selection to match the author's final view state.
$ws.Range("A83").Select()
